$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")

# --- ALC ---
# row 82
$ws1.Range("H82").Value = 1981.4286
$ws1.Range("I82").Value = 1981.4286
$ws1.Range("K82").Value = 5944.2858
$ws1.Range("M82").Value = -5538.2858
# row 85
$ws1.Range("H85").Value = 1981.4286
$ws1.Range("I85").Value = 1981.4286
$ws1.Range("K85").Value = 5944.2858
$ws1.Range("M85").Value = -4540.2858
# row 86
$ws1.Range("H86").Value = 4347.9
$ws1.Range("I86").Value = 2996.5
$ws1.Range("K86").Value = 2996.5
$ws1.Range("M86").Value = -1873.5
# row 89
$ws1.Range("H89").Value = 4347.9
$ws1.Range("I89").Value = 2996.5
$ws1.Range("K89").Value = 14982.5
$ws1.Range("M89").Value = -9366.5
# row 98
$ws1.Range("H98").Value = 2716.125
$ws1.Range("I98").Value = 2747
$ws1.Range("K98").Value = 2747
$ws1.Range("M98").Value = -1249
# row 107
$ws1.Range("H107").Value = 1384.5
$ws1.Range("I107").Value = 1735.0714
$ws1.Range("K107").Value = 1735.0714
$ws1.Range("M107").Value = 184.9286
# row 111
$ws1.Range("H111").Value = 3550.182
$ws1.Range("I111").Value = 416.5
$ws1.Range("K111").Value = 1249.5
$ws1.Range("M111").Value = 1817.5
# row 122
$ws1.Range("H122").Value = 2716.125
$ws1.Range("I122").Value = 2747
$ws1.Range("K122").Value = 8241
$ws1.Range("M122").Value = -5791
# row 132
$ws1.Range("H132").Value = 14173.116
$ws1.Range("I132").Value = 1847.7097
$ws1.Range("K132").Value = 5543.1291
$ws1.Range("M132").Value = -3013.1291
# row 138
$ws1.Range("H138").Value = 3473.6667
$ws1.Range("I138").Value = 2267
$ws1.Range("J138").Value = 3818.4285
$ws1.Range("K138").Value = 6801
$ws1.Range("L138").Value = 11455.2855
$ws1.Range("M138").Value = -1661
$ws1.Range("N138").Value = -21735.2855

# --- ARM ---
# row 32
$ws2.Range("H32").Value = 5087.274
$ws2.Range("I32").Value = 5177.443
$ws2.Range("K32").Value = 5177.443
$ws2.Range("M32").Value = -4890.443
# row 88
$ws2.Range("H88").Value = 1140.3125
$ws2.Range("I88").Value = 1650.6666
$ws2.Range("J88").Value = 834.1
$ws2.Range("K88").Value = 1650.6666
$ws2.Range("L88").Value = 834.1
$ws2.Range("M88").Value = -1244.6666
$ws2.Range("N88").Value = -1646.1
# row 91
$ws2.Range("H91").Value = 1140.3125
$ws2.Range("I91").Value = 1650.6666
$ws2.Range("J91").Value = 834.1
$ws2.Range("K91").Value = 1650.6666
$ws2.Range("L91").Value = 834.1
$ws2.Range("M91").Value = -246.6666
$ws2.Range("N91").Value = -3642.1
# row 102
$ws2.Range("H102").Value = 3865.3333
$ws2.Range("I102").Value = 3865.3333
$ws2.Range("K102").Value = 3865.3333
$ws2.Range("M102").Value = -2243.3333
# row 122
$ws2.Range("H122").Value = 4380.5527
$ws2.Range("I122").Value = 3827.0908
$ws2.Range("K122").Value = 11481.2724
$ws2.Range("M122").Value = -9031.2724

# --- BSM ---
# row 94
$ws3.Range("H94").Value = 4297.522
$ws3.Range("I94").Value = 3378.75
$ws3.Range("K94").Value = 3378.75
$ws3.Range("M94").Value = -2927.75
# row 107
$ws3.Range("H107").Value = 3092.5715
$ws3.Range("I107").Value = 2789.3713
$ws3.Range("K107").Value = 2789.3713
$ws3.Range("M107").Value = -869.3712999999998
# row 133
$ws3.Range("H133").Value = 75000
$ws3.Range("J133").Value = 75000
$ws3.Range("L133").Value = 75000
$ws3.Range("N133").Value = -85120

# --- CRP ---
# row 22
$ws4.Range("H22").Value = 456
$ws4.Range("I22").Value = 525.6667
$ws4.Range("J22").Value = 316.66666
$ws4.Range("K22").Value = 525.6667
$ws4.Range("L22").Value = 316.66666
$ws4.Range("M22").Value = -175.6667
$ws4.Range("N22").Value = -1016.66666
# row 31
$ws4.Range("H31").Value = 1431.5143
$ws4.Range("I31").Value = 1303.4667
$ws4.Range("J31").Value = 2199.8
$ws4.Range("K31").Value = 1303.4667
$ws4.Range("L31").Value = 2199.8
$ws4.Range("M31").Value = -1008.4667
$ws4.Range("N31").Value = -2789.8
# row 34
$ws4.Range("H34").Value = 1431.5143
$ws4.Range("I34").Value = 1303.4667
$ws4.Range("J34").Value = 2199.8
$ws4.Range("K34").Value = 1303.4667
$ws4.Range("L34").Value = 2199.8
$ws4.Range("M34").Value = -1101.4667
$ws4.Range("N34").Value = -2603.8
# row 62
$ws4.Range("H62").Value = 50002904
$ws4.Range("I62").Value = 3334
$ws4.Range("J62").Value = 166668580
$ws4.Range("K62").Value = 3334
$ws4.Range("L62").Value = 166668580
$ws4.Range("M62").Value = -2710
$ws4.Range("N62").Value = -166669828
# row 65
$ws4.Range("H65").Value = 50002904
$ws4.Range("I65").Value = 3334
$ws4.Range("J65").Value = 166668580
$ws4.Range("K65").Value = 16670
$ws4.Range("L65").Value = 833342900
$ws4.Range("M65").Value = -13550
$ws4.Range("N65").Value = -833349140
# row 99
$ws4.Range("H99").Value = 3451.182
$ws4.Range("I99").Value = 3507.0908
$ws4.Range("K99").Value = 3507.0908
$ws4.Range("M99").Value = -2009.0908
# row 105
$ws4.Range("H105").Value = 525
$ws4.Range("I105").Value = 493.33334
$ws4.Range("K105").Value = 493.33334
$ws4.Range("M105").Value = 1253.66666
# row 126
$ws4.Range("H126").Value = 3451.182
$ws4.Range("I126").Value = 3507.0908
$ws4.Range("K126").Value = 10521.2724
$ws4.Range("M126").Value = -8051.2724
# row 132
$ws4.Range("H132").Value = 2543.3823
$ws4.Range("I132").Value = 1862.2858
$ws4.Range("K132").Value = 5586.857400000001
$ws4.Range("M132").Value = -3056.857400000001

# --- CUL ---
# row 92
$ws5.Range("H92").Value = 0
$ws5.Range("J92").Value = 0
$ws5.Range("L92").Value = 0
$ws5.Range("N92").ClearContents()
# row 94
$ws5.Range("H94").Value = 3831
$ws5.Range("J94").Value = 3831
$ws5.Range("L94").Value = 11493
$ws5.Range("N94").Value = -12845
# row 95
$ws5.Range("H95").Value = 4950
$ws5.Range("J95").Value = 4950
$ws5.Range("L95").Value = 14850
$ws5.Range("N95").Value = -18968
# row 107
$ws5.Range("H107").Value = 2098
$ws5.Range("I107").Value = 2387.2222
$ws5.Range("J107").Value = 1967.85
$ws5.Range("K107").Value = 7161.6666
$ws5.Range("L107").Value = 5903.549999999999
$ws5.Range("M107").Value = -5241.6666
$ws5.Range("N107").Value = -9743.549999999999

# --- GSM ---
# row 33
$ws6.Range("H33").Value = 25499
$ws6.Range("J33").Value = 25499
$ws6.Range("L33").Value = 25499
$ws6.Range("N33").Value = -26003
# row 44
$ws6.Range("H44").Value = 23312.666
$ws6.Range("J44").Value = 25999
$ws6.Range("L44").Value = 25999
$ws6.Range("N44").Value = -27191
# row 47
$ws6.Range("H47").Value = 19999.5
$ws6.Range("J47").Value = 19999.5
$ws6.Range("L47").Value = 19999.5
$ws6.Range("N47").Value = -21135.5
# row 80
$ws6.Range("H80").Value = 49026.555
$ws6.Range("I80").Value = 71605.766
$ws6.Range("J80").Value = 10641.9
$ws6.Range("K80").Value = 71605.766
$ws6.Range("L80").Value = 10641.9
$ws6.Range("M80").Value = -70607.766
$ws6.Range("N80").Value = -12637.9
# row 83
$ws6.Range("H83").Value = 49026.555
$ws6.Range("I83").Value = 71605.766
$ws6.Range("J83").Value = 10641.9
$ws6.Range("K83").Value = 358028.83
$ws6.Range("L83").Value = 53209.5
$ws6.Range("M83").Value = -353036.83
$ws6.Range("N83").Value = -63193.5

# --- LTW ---
# row 82
$ws7.Range("H82").Value = 66668692
$ws7.Range("I82").Value = 111113550
$ws7.Range("J82").Value = 1404.1666
$ws7.Range("K82").Value = 111113550
$ws7.Range("L82").Value = 1404.1666
$ws7.Range("M82").Value = -111113189
$ws7.Range("N82").Value = -2126.1666
# row 85
$ws7.Range("H85").Value = 66668692
$ws7.Range("I85").Value = 111113550
$ws7.Range("J85").Value = 1404.1666
$ws7.Range("K85").Value = 111113550
$ws7.Range("L85").Value = 1404.1666
$ws7.Range("M85").Value = -111112302
$ws7.Range("N85").Value = -3900.1666
# row 122
$ws7.Range("H122").Value = 5981.5
$ws7.Range("I122").Value = 5571.1816
$ws7.Range("J122").Value = 6626.2856
$ws7.Range("K122").Value = 16713.5448
$ws7.Range("L122").Value = 19878.8568
$ws7.Range("M122").Value = -14263.5448
$ws7.Range("N122").Value = -24778.8568
# row 132
$ws7.Range("H132").Value = 2436.0952
$ws7.Range("I132").Value = 2057.9
$ws7.Range("K132").Value = 6173.700000000001
$ws7.Range("M132").Value = -3643.700000000001

# --- WVR ---
# row 81
$ws8.Range("H81").Value = 2114.5
$ws8.Range("I81").Value = 1586.3334
$ws8.Range("K81").Value = 3172.6668
$ws8.Range("M81").Value = -2111.6668
# row 84
$ws8.Range("H84").Value = 2114.5
$ws8.Range("I84").Value = 1586.3334
$ws8.Range("K84").Value = 15863.334
$ws8.Range("M84").Value = -10559.334
# row 96
$ws8.Range("H96").Value = 45467.332
$ws8.Range("I96").Value = 86792.336
$ws8.Range("K96").Value = 86792.336
$ws8.Range("M96").Value = -85419.336

Write-Output "Applied all updates"